$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edit -----------------------------------------------------------
# The author selected the block of dates/prices from 2017-10-14 through
# 2017-12-12 (rows 288-347, columns A:B) and cleared the cell contents.
# ClearContents removes the values but keeps the existing cell formatting
# (column A keeps its date number format, style index 1).
$ws.Range("A288:B347").ClearContents()

# --- View / selection edit ------------------------------------------------
# Scroll the sheet so row 274 is at the top and select A293, then restore
# the workbook window position recorded the last time it was saved.
$ws.Activate()
try { $excel.ActiveWindow.ScrollColumn = 1 } catch { }
try { $excel.ActiveWindow.ScrollRow = 274 } catch { }
$ws.Range("A293").Select()

try { $excel.ActiveWindow.Left = 1280 } catch { }
try { $excel.ActiveWindow.Top = 1200 } catch { }
try { $excel.Left = 1280 } catch { }
try { $excel.Top = 1200 } catch { }
